$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Date" (the EquipmentInventory.dateReceived type) -> "Integer"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Date", $true, $true, $false, $false, $false, `
    $true, 1, $false, "Integer", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Intro paragraph: add a sentence about Challenges gradually
#    including less guidance.
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "the intended SQL commands for each task are included.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "the intended SQL commands for each task are included. Challenges will gradually include less guidance and more complex queries.", `
    2) | Out-Null

# ------------------------------------------------------------------
# 3) "the New Employee Orientation. " -> "...packet. "
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "the New Employee Orientation. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "the New Employee Orientation packet. ", `
    2) | Out-Null

# ------------------------------------------------------------------
# 4) "...boot it up inside a virtual machine so you can investigate."
#    -> "...investigate later."
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "boot it up inside a virtual machine so you can investigate.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "boot it up inside a virtual machine so you can investigate later.", `
    2) | Out-Null

# ------------------------------------------------------------------
# 5) New "Day Two: Intruder" section appended at the end of the doc.
# ------------------------------------------------------------------

# -- Heading2: "Day Two: Intruder"
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Style = "Heading 2"
$headingPara.Range.InsertAfter("Day Two: Intruder")

# -- Narrative paragraph
$headingPara2 = $d.Paragraphs.Last
$headingPara2.Range.InsertParagraphAfter()
$introPara = $d.Paragraphs.Last
$introPara.Range.InsertAfter( `
    "You find that this computer was given to an employee named " + `
    [char]0x201C + "Cal Irris" + [char]0x201D + `
    " after it was already infected. This occurred sometime between the dates 78-554-210-6 and 82-974-000-2. However, that name does not match any employees in the database. You want to find all the IDs of all employees that were given computers between those dates, so you can narrow down the owner" + `
    [char]0x2019 + "s real name. ")

# -- "Potential Solution" bullet
$introPara2 = $d.Paragraphs.Last
$introPara2.Range.InsertParagraphAfter()
$solPara = $d.Paragraphs.Last
$solPara.Style = "List Paragraph"
$solPara.Range.ListFormat.ApplyBulletDefault()
$solRange = $solPara.Range
$solRange.Collapse(0)
$solRange.InsertAfter("Potential Solution")
$solRange.Font.Italic = 1
$solRange.Font.ItalicBi = 1

# -- SQL lines, each indented 0.5"
$solPara2 = $d.Paragraphs.Last
$solPara2.Range.InsertParagraphAfter()
$sql1 = $d.Paragraphs.Last
$sql1.Range.InsertAfter("SELECT DISTINCT employeeID")
$sql1.LeftIndent = 36

$sql1b = $d.Paragraphs.Last
$sql1b.Range.InsertParagraphAfter()
$sql2 = $d.Paragraphs.Last
$sql2.Range.InsertAfter("FROM EquipmentInventory")
$sql2.LeftIndent = 36

$sql2b = $d.Paragraphs.Last
$sql2b.Range.InsertParagraphAfter()
$sql3 = $d.Paragraphs.Last
$sql3.Range.InsertAfter("WHERE category = " + [char]0x201C + "Computer" + [char]0x201D)
$sql3.LeftIndent = 36

$sql3b = $d.Paragraphs.Last
$sql3b.Range.InsertParagraphAfter()
$sql4 = $d.Paragraphs.Last
$sql4.Range.InsertAfter("GROUP BY employeeID")
$sql4.LeftIndent = 36

$sql4b = $d.Paragraphs.Last
$sql4b.Range.InsertParagraphAfter()
$sql5 = $d.Paragraphs.Last
$sql5.Range.InsertAfter("HAVING dateReceived BETWEEN 785542106 AND 829740002")
$sql5.LeftIndent = 36

# -- two trailing empty paragraphs
$sql5b = $d.Paragraphs.Last
$sql5b.Range.InsertParagraphAfter()
$sql5b.Range.InsertParagraphAfter()

Write-Output "done"
